# Applies the Sargatanas_Profits.xlsx profit-table refresh (scheduled runner update).
# For each touched row, H/I/J/K/L/M/N are rewritten to the new computed values;
# some cells are newly populated (were previously absent) and some are cleared
# (value removed entirely) to mirror the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H19").Value = 3570.7144
$ws.Range("I19").Value = 3365.3333
$ws.Range("J19").Value = 3724.75
$ws.Range("K19").Value = 3365.3333
$ws.Range("L19").Value = 3724.75
$ws.Range("M19").Value = -3190.3333
$ws.Range("N19").Value = -4074.75

$ws.Range("H43").Value = 342668.5
$ws.Range("I43").Value = 616.25
$ws.Range("K43").Value = 616.25
$ws.Range("M43").Value = -547.25

$ws.Range("H48").Value = 746.25
$ws.Range("J48").Value = 290
$ws.Range("L48").Value = 870
$ws.Range("N48").Value = -1454

$ws.Range("H56").Value = 746.25
$ws.Range("J56").Value = 290
$ws.Range("L56").Value = 870
$ws.Range("N56").Value = -1938

$ws.Range("H62").Value = 76951580
$ws.Range("I62").Value = 200000740
$ws.Range("K62").Value = 200000740
$ws.Range("M62").Value = -200000116

$ws.Range("H65").Value = 76951580
$ws.Range("I65").Value = 200000740
$ws.Range("K65").Value = 1000003700
$ws.Range("M65").Value = -1000000580

$ws.Range("H76").Value = 100000000
$ws.Range("I76").Value = 100000000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 100000000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -99999685
$ws.Range("N76").Value = ""

$ws.Range("H79").Value = 100000000
$ws.Range("I79").Value = 100000000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 100000000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -99998908
$ws.Range("N79").Value = ""

$ws.Range("H86").Value = 65591124
$ws.Range("I86").Value = 93753780
$ws.Range("K86").Value = 93753780
$ws.Range("M86").Value = -93752657

$ws.Range("H87").Value = 60000
$ws.Range("J87").Value = 60000
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62496

$ws.Range("H89").Value = 65591124
$ws.Range("I89").Value = 93753780
$ws.Range("K89").Value = 468768900
$ws.Range("M89").Value = -468763284

$ws.Range("H90").Value = 60000
$ws.Range("J90").Value = 60000
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -192480

$ws.Range("H94").Value = 5707.4443
$ws.Range("I94").Value = 195.28572
$ws.Range("K94").Value = 195.28572
$ws.Range("M94").Value = 255.71428

$ws.Range("H103").Value = 1188.7059
$ws.Range("J103").Value = 1388.0769
$ws.Range("L103").Value = 4164.2307
$ws.Range("N103").Value = -5336.2307

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").Value = ""

$ws.Range("H107").Value = 22116438
$ws.Range("I107").Value = 10715400
$ws.Range("J107").Value = 70000800
$ws.Range("K107").Value = 10715400
$ws.Range("L107").Value = 70000800
$ws.Range("M107").Value = -10713480
$ws.Range("N107").Value = -70004640

$ws.Range("H132").Value = 1998.65
$ws.Range("I132").Value = 1998.2632
$ws.Range("K132").Value = 5994.7896
$ws.Range("M132").Value = -3464.7896

$ws.Range("H138").Value = 4189.9185
$ws.Range("I138").Value = 1912.6666
$ws.Range("J138").Value = 5512.1934
$ws.Range("K138").Value = 5737.9998
$ws.Range("L138").Value = 16536.5802
$ws.Range("M138").Value = -597.9997999999996
$ws.Range("N138").Value = -26816.5802

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 3008.48
$ws.Range("I32").Value = 2942.6345
$ws.Range("J32").Value = 3883.2856
$ws.Range("K32").Value = 2942.6345
$ws.Range("L32").Value = 3883.2856
$ws.Range("M32").Value = -2655.6345
$ws.Range("N32").Value = -4457.2856

$ws.Range("H122").Value = 13882.556
$ws.Range("I122").Value = 15492
$ws.Range("K122").Value = 46476
$ws.Range("M122").Value = -44026

$ws.Range("H132").Value = 2846.4521
$ws.Range("I132").Value = 1150.2115
$ws.Range("J132").Value = 7046.6665
$ws.Range("K132").Value = 3450.6345
$ws.Range("L132").Value = 21139.9995
$ws.Range("M132").Value = -920.6344999999997
$ws.Range("N132").Value = -26199.9995

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H20").Value = 10418702
$ws.Range("I20").Value = 33336106
$ws.Range("J20").Value = 1700.6364
$ws.Range("K20").Value = 33336106
$ws.Range("L20").Value = 1700.6364
$ws.Range("M20").Value = -33335859
$ws.Range("N20").Value = -2194.6364

$ws.Range("H80").Value = 41667010
$ws.Range("J80").Value = 297.14285
$ws.Range("L80").Value = 297.14285
$ws.Range("N80").Value = -2293.14285

$ws.Range("H83").Value = 41667010
$ws.Range("J83").Value = 297.14285
$ws.Range("L83").Value = 1485.71425
$ws.Range("N83").Value = -11469.71425

$ws.Range("H105").Value = 1999.4
$ws.Range("I105").Value = 1777.1111
$ws.Range("K105").Value = 1777.1111
$ws.Range("M105").Value = -30.11110000000008

$ws.Range("H107").Value = 37503412
$ws.Range("J107").Value = 5534.2
$ws.Range("L107").Value = 5534.2
$ws.Range("N107").Value = -9374.200000000001

$ws.Range("H134").Value = 4155.082
$ws.Range("I134").Value = 1094.5834
$ws.Range("J134").Value = 8562.200000000001
$ws.Range("K134").Value = 3283.7502
$ws.Range("L134").Value = 25686.6
$ws.Range("M134").Value = -748.7501999999999
$ws.Range("N134").Value = -30756.6

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H132").Value = 4360.28
$ws.Range("I132").Value = 1811.3243
$ws.Range("J132").Value = 11615
$ws.Range("K132").Value = 5433.9729
$ws.Range("L132").Value = 34845
$ws.Range("M132").Value = -2903.9729
$ws.Range("N132").Value = -39905

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H56").Value = 7093.8
$ws.Range("I56").Value = 7093.8
$ws.Range("K56").Value = 7093.8
$ws.Range("M56").Value = -6563.8

$ws.Range("H62").Value = 2999.6667
$ws.Range("J62").Value = 2999
$ws.Range("L62").Value = 8997
$ws.Range("N62").Value = -10369

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = ""

$ws.Range("H64").Value = 16668374
$ws.Range("I64").Value = 1415
$ws.Range("J64").Value = 33335332
$ws.Range("K64").Value = 4245
$ws.Range("L64").Value = 100005996
$ws.Range("M64").Value = -3975
$ws.Range("N64").Value = -100006536

$ws.Range("H65").Value = 2999.6667
$ws.Range("J65").Value = 2999
$ws.Range("L65").Value = 26991
$ws.Range("N65").Value = -33855

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = ""

$ws.Range("H67").Value = 16668374
$ws.Range("I67").Value = 1415
$ws.Range("J67").Value = 33335332
$ws.Range("K67").Value = 4245
$ws.Range("L67").Value = 100005996
$ws.Range("M67").Value = -3309
$ws.Range("N67").Value = -100007868

$ws.Range("H70").Value = 4000
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -2685
$ws.Range("N70").Value = -15630

$ws.Range("H73").Value = 4000
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -1908
$ws.Range("N73").Value = -17184

$ws.Range("H126").Value = 2933
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = ""

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H2").Value = 43.411766
$ws.Range("I2").Value = 33.6
$ws.Range("K2").Value = 33.6
$ws.Range("M2").Value = 79.40000000000001

$ws.Range("H133").Value = 77780
$ws.Range("J133").Value = 77780
$ws.Range("L133").Value = 77780
$ws.Range("N133").Value = -87900

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").Value = ""

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").Value = ""

$ws.Range("H68").Value = 2702.8235
$ws.Range("I68").Value = 1579
$ws.Range("K68").Value = 1579
$ws.Range("M68").Value = -830

$ws.Range("H71").Value = 2702.8235
$ws.Range("I71").Value = 1579
$ws.Range("K71").Value = 7895
$ws.Range("M71").Value = -4151

$ws.Range("H87").Value = 56500
$ws.Range("J87").Value = 56500
$ws.Range("L87").Value = 56500
$ws.Range("N87").Value = -58746

$ws.Range("H90").Value = 56500
$ws.Range("J90").Value = 56500
$ws.Range("L90").Value = 169500
$ws.Range("N90").Value = -180732

$ws.Range("H100").Value = 5107.3335
$ws.Range("I100").Value = 3897.8333
$ws.Range("K100").Value = 3897.8333
$ws.Range("M100").Value = -3356.8333

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H136").Value = 16687002
$ws.Range("I136").Value = 25001110
$ws.Range("K136").Value = 75003330
$ws.Range("M136").Value = -75000780
